# Fruta / hortaliza, semanal
# Insert a new weekly record as row 153 (Femacal de La Calera, Coquimbo,
# Arándano (blue)), pushing the existing rows 153-220 down to 154-221.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new blank row at position 153, shifting rows 153:220 down to 154:221
$ws.Rows(153).Insert()

# Populate the newly inserted row with the new weekly observation
$ws.Range("A153").Value = 3
$ws.Range("B153").Value = 'Femacal de La Calera'
$ws.Range("C153").Value = 'Coquimbo'
$ws.Range("D153").Value = 44813
$ws.Range("E153").Value = 5
$ws.Range("F153").Value = 'Fruta'
$ws.Range("G153").Value = 100101
$ws.Range("H153").Value = 'Berries'
$ws.Range("I153").Value = 100101001
$ws.Range("J153").Value = 'Arándano (blue)'
$ws.Range("K153").Value = 'Sin especificar'
$ws.Range("L153").Value = 'Primera'
$ws.Range("M153").Value = 45
$ws.Range("N153").Value = 12000
$ws.Range("O153").Value = 12000
$ws.Range("P153").Value = 12000
$ws.Range("Q153").Value = '$/bandeja 12 canastillos 125 gramos'
$ws.Range("R153").Value = 'Provincia de Limarí'
$ws.Range("S153").Value = 8000
$ws.Range("T153").Value = 1.5
